$d = $word.ActiveDocument

# Clear all existing content of the single paragraph (this also removes the
# "_GoBack" bookmark, which we re-create below at its correct new location,
# right after "working with angularJS" in the rewritten 5th paragraph).
$d.Content.Delete()

# Rebuild the whole body as one text assignment: backtick-r ("`r") marks a
# paragraph break, producing the full 27-paragraph structure (14 paragraphs
# of text interleaved with 13 blank paragraphs) described by the diff.
$r = $d.Paragraphs(1).Range
$r.Text = "Hello everyone! My name is Malachi Gray and I’m a graduate of the rails engineering course at The Iron Yard this winter.  `r`rPrior to attending the Iron Yard, I was a project analyst at a contractor working for the government.  After 6 years I decided I wanted to do something more technical and I also wanted to find a career that allowed me to continually learn and evolve.  In December I left my job and with my wife, moved to Charleston to attend the Iron Yard and begin my new career in the awesome world of programming.`r`rAs Julie mentioned, we built the Charleston Basketlist app which allows users to select activities they would like to complete around Charleston, like a personal bucketlist.`r`rThis application is an integrated Rails and AngularJs application.  `r`rMy core responsibility was developing the rails side functionality but Julie also gave me the latitude to dive right in to working with angularJS.  `r`rOn the rails side, we have multiple models, including join tables and nesting.  `r`rIn addition, we maintain user authentication using devise and authorizations using cancancan on the server side.  `r`rModels and controllers were tested in RSPEC.`r`rTwo of the core front end items I took the lead on were user comments and image uploads to AWS S3 buckets.  `r`rIn the comments section, users are able to create and see comments and then are given the opportunity to edit or delete only if they made that comment.`r`rFor the AWS S3 bucket uploads, we were able to incorporate the ngUpload directive, developed here in Charleston, to send an image to the rails server which uses paperclip to send it to AWS S3.`r`rThe most challenging part of the project was learning how rails and angular interacted and integrated within the Rails pipeline.  We found out early on there would be quite a bit of learning involved with this project but we were able to work together and make things work. `r`rI know both Julie and I learned a lot on this project and both enjoyed it a lot. `r`rThanks for your time!  We would love to show you our app at the table."

Write-Output "Paragraphs.Count=$($d.Paragraphs.Count)"

# Re-create the _GoBack bookmark exactly where it now belongs: immediately
# after "working with angularJS" (before the closing ".  " of that
# paragraph), matching its position in the target markup.
$bmRange = $d.Content
$bmRange.Find.Execute("working with angularJS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
Write-Output "Bookmark added at $($bmRange.Start)"
